$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (45243 = 2023-11-13) for every
# data row (2-83). The update bumps that date forward by one day (45244 =
# 2023-11-14) for all of them.
$ws.Range("C2:C83").Value = 45244
